$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column AT (row 1) - next date label after "04-ago"
$ws.Cells.Item(1, 46).Value = "07-ago"

# New data values for column AT, rows 2-18
$atValues = @{
    2  = 0
    3  = 14.103146384608888
    4  = 19.175891190222764
    5  = 17.72017533329516
    6  = 0
    7  = 14.672089134327582
    8  = 15.75586833200383
    9  = 11.856526897434366
    10 = 11.368491737273116
    11 = 14.49408607939103
    12 = 0
    13 = 5.8945829087612465
    14 = 0
    15 = 0
    16 = 12.567299766239854
    17 = 0
    18 = 0
}

foreach ($row in $atValues.Keys) {
    $ws.Cells.Item($row, 46).Value = $atValues[$row]
}

# Match the final cell selection left behind in the saved workbook
$null = $ws.Range("AV5").Select()
